$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): columns shifted / renamed ---------------------
# C1: "Dauer" -> "Schaden"
$ws.Range("C1").Value = "Schaden"
# D1: "Kosten" -> "QS"
$ws.Range("D1").Value = "QS"
# E1: "pro Stunde" -> "Zauberdauer"
$ws.Range("E1").Value = "Zauberdauer"
# F1: "Reichweite" -> "Kosten"
$ws.Range("F1").Value = "Kosten"
# G1: "Wirkdauer" -> "pro Stunde"
$ws.Range("G1").Value = "pro Stunde"

# --- Row 2: Bannbaladin --------------------------------------------------
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "4 Aktion(en)"
$ws.Range("F2").Value = "8 AsP"
$ws.Range("G2").Value = "0 AsP"

# --- Row 3: Flim Flam -----------------------------------------------------
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "1 Aktion(en)"
$ws.Range("F3").Value = "2 AsP"
$ws.Range("G3").Value = "1 AsP"

# --- Row 4: Ignifaxius (new spell method: magischer angriff) -------------
$ws.Range("C4").Value = "2 w6"
$ws.Range("D4").Value = "x2"
$ws.Range("E4").Value = "2 Aktion(en)"
$ws.Range("F4").Value = "8 AsP"
$ws.Range("G4").Value = "0 AsP"

# --- Formatting: center-align the whole used range ------------------------
$ws.Range("A1:G4").HorizontalAlignment = -4108

# --- Selection moves to G9 -------------------------------------------------
$ws.Range("G9").Select()
